$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item("004江润洲")
$ws4.Range("B6").Value = "'13"
